$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = 'Kenya'
$ws.Range("B2").Value = "'4"
$ws.Range("C2").Value = 'Food Processing'
$ws.Range("D2").Value = 0.02275
$ws.Range("G2").Value = -0.001861407278649935
$ws.Range("H2").Value = -0.001861407278649935
$ws.Range("I2").Value = -0.03318800632395387
$ws.Range("J2").Value = -0.03318800632395387
$ws.Range("K2").Value = -0.552
$ws.Range("L2").Value = -0.007555123660402667
$ws.Range("M2").Value = 4.768
$ws.Range("N2").Value = 0.06727811485819105
$ws.Range("O2").Value = -8.637681159420289
$ws.Range("P2").Value = 4.768
$ws.Range("Q2").Value = 0.06727811485819105
$ws.Range("R2").Value = -8.637681159420289
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 5.14
$ws.Range("V2").Value = 0.07252716241004656
$ws.Range("X2").Value = 0.08038484172650837
$ws.Range("Z2").Value = 0.7367146996092802
$ws.Range("AB2").Value = 0.08038484172650837
$ws.Range("AD2").Value = 4.34
$ws.Range("AE2").Value = 0.5940765302352066
$ws.Range("AF2").Value = 4.934076530235206
$ws.Range("AG2").Value = -0.2059234697647936
$ws.Range("AH2").Value = 0.06508985738078613
$ws.Range("AI2").Value = 0.03868829935084415
$ws.Range("AJ2").Value = -0.002914118175402527
$ws.Range("AK2").Value = -0.001682462710635543
$ws.Range("AL2").Value = 0.206
$ws.Range("AM2").Value = -0.875
$ws.Range("AN2").Value = 1.000691722388748
$ws.Range("AO2").Value = -12.01456310679612
$ws.Range("AP2").Value = -0.04748062480165866
$ws.Range("AQ2").Value = 2.828571428571429
$ws.Range("AA2").ClearContents()
$ws.Range("AC2").ClearContents()

# --- Row 3 ---
$ws.Range("A3").Value = 'Kenya'
$ws.Range("B3").Value = 'Sasini PLC (NASE:SASN)'
$ws.Range("C3").Value = 'Food Processing'
$ws.Range("D3").Value = 0.0524
$ws.Range("G3").Value = -0.1168711656441718
$ws.Range("H3").Value = -0.1168711656441718
$ws.Range("I3").Value = -0.03281642043089084
$ws.Range("J3").Value = -0.03281642043089084
$ws.Range("K3").Value = -1.34
$ws.Range("L3").Value = -0.04110429447852761
$ws.Range("M3").Value = 0.751
$ws.Range("N3").Value = 0.01845208845208845
$ws.Range("O3").Value = -0.5604477611940298
$ws.Range("P3").Value = 0.751
$ws.Range("Q3").Value = 0.01845208845208845
$ws.Range("R3").Value = -0.5604477611940298
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 5.14
$ws.Range("V3").Value = 0.1262899262899263
$ws.Range("W3").Value = -0.01283524904214559
$ws.Range("X3").Value = 0.08641719146736276
$ws.Range("Y3").Value = -0.09925244050950836
$ws.Range("Z3").Value = 0.3287149337867666
$ws.Range("AA3").Value = -0.01078724746905898
$ws.Range("AB3").Value = 0.08302338403178902
$ws.Range("AC3").Value = -0.09381063150084801
$ws.Range("AD3").Value = 4.34
$ws.Range("AE3").Value = 0.5940765302352066
$ws.Range("AF3").Value = 4.934076530235206
$ws.Range("AG3").Value = -0.2059234697647936
$ws.Range("AH3").Value = 0.1081226334659385
$ws.Range("AI3").Value = 0.03868829935084415
$ws.Range("AJ3").Value = -0.005085273882243968
$ws.Range("AK3").Value = -0.001682462710635543
$ws.Range("AL3").Value = 0.206
$ws.Range("AM3").Value = -0.01000000000000001
$ws.Range("AN3").Value = 6.546003016591252
$ws.Range("AO3").Value = -5.436893203883496
$ws.Range("AP3").Value = -0.3105934687251788
$ws.Range("AQ3").Value = 111.9999999999999

# --- Row 4 ---
$ws.Range("A4").Value = 'Kenya'
$ws.Range("B4").Value = 'Williamson Tea Kenya Plc (NASE:WTK)'
$ws.Range("C4").Value = 'Food Processing'
$ws.Range("D4").Value = 0.0345
$ws.Range("G4").Value = 0.1458904109589041
$ws.Range("H4").Value = 0.1458904109589041
$ws.Range("I4").Value = -0.01092465753424658
$ws.Range("J4").Value = -0.01092465753424658
$ws.Range("K4").Value = 1.27
$ws.Range("L4").Value = 0.04349315068493151
$ws.Range("M4").Value = 3.28
$ws.Range("N4").Value = 0.1569377990430622
$ws.Range("O4").Value = 2.582677165354331
$ws.Range("P4").Value = 3.28
$ws.Range("Q4").Value = 0.1569377990430622
$ws.Range("R4").Value = 2.582677165354331
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("X4").Value = 0.08038484172650837
$ws.Range("AB4").Value = 0.08038484172650837
$ws.Range("AD4").Value = 0
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AJ4").Value = 0
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = -0.574
$ws.Range("AN4").Value = 0
$ws.Range("AP4").Value = 0
$ws.Range("AQ4").Value = 0.5557491289198607
$ws.Range("Z4").ClearContents()
$ws.Range("AA4").ClearContents()
$ws.Range("AC4").ClearContents()
$ws.Range("AI4").ClearContents()
$ws.Range("AK4").ClearContents()
$ws.Range("AO4").ClearContents()

# --- Row 5 ---
$ws.Range("A5").Value = 'Kenya'
$ws.Range("B5").Value = 'Kapchorua Tea Kenya Plc (NASE:KAPC)'
$ws.Range("C5").Value = 'Food Processing'
$ws.Range("D5").Value = 0.011
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = -0.02787037037037037
$ws.Range("J5").Value = -0.02787037037037037
$ws.Range("K5").Value = 0.185
$ws.Range("L5").Value = 0.01712962962962963
$ws.Range("M5").Value = 0.737
$ws.Range("N5").Value = 0.1318425760286225
$ws.Range("O5").Value = 3.983783783783784
$ws.Range("P5").Value = 0.737
$ws.Range("Q5").Value = 0.1318425760286225
$ws.Range("R5").Value = 3.983783783783784
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 0
$ws.Range("V5").Value = 0
$ws.Range("X5").Value = 0.08038484172650837
$ws.Range("AB5").Value = 0.08038484172650837
$ws.Range("AD5").Value = 0
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 0
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AJ5").Value = 0
$ws.Range("AL5").Value = 0
$ws.Range("AM5").Value = -0.291
$ws.Range("AQ5").Value = 1.034364261168385

# --- Row 6 ---
$ws.Range("A6").Value = 'Kenya'
$ws.Range("B6").Value = 'Eaagads Limited (NASE:EGAD)'
$ws.Range("C6").Value = 'Food Processing'
$ws.Range("D6").Value = -0.137
$ws.Range("G6").Value = -1.265658747300216
$ws.Range("H6").Value = -1.265658747300216
$ws.Range("I6").Value = -1.587473002159827
$ws.Range("J6").Value = -1.587473002159827
$ws.Range("K6").Value = -0.667
$ws.Range("L6").Value = -1.44060475161987
$ws.Range("M6").Value = -0
$ws.Range("N6").Value = -0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = -0
$ws.Range("Q6").Value = -0
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("U6").Value = 0
$ws.Range("V6").Value = 0
$ws.Range("X6").Value = 0.08038484172650837
$ws.Range("AB6").Value = 0.08038484172650837
$ws.Range("AD6").Value = 0
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 0
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AJ6").Value = 0
$ws.Range("AL6").Value = 0
$ws.Range("AM6").Value = 0
$ws.Range("AN6").Value = -0
$ws.Range("AP6").Value = -0
